# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet's column G (header "K") held per-game strikeout counts that were
# previously computed a different way ("Strike#"); this regenerates those
# values to the corrected K counts for each row (game).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    7  = 3
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 3
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 4
    22 = 1
    23 = 3
    24 = 2
    25 = 1
    26 = 1
    27 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
